$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.497786283493042
$ws.Range("B1").Value = 3.811347961425781
$ws.Range("C1").Value = 2.741596460342407
$ws.Range("D1").Value = 0.7987833619117737
$ws.Range("E1").Value = 1.073408007621765
